$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "test 1" (B2) and "test 1 fr" (D2) are no longer translated for
# this locale pass -- the cells are dropped entirely (content + formatting).
$ws.Range("B2").Clear()
$ws.Range("D2").Clear()

# Row 2: de_DE (H2) translation was pulled, leaving the cell blank; el (I2)
# gets its real translation. Both lose the old bold/header-ish formatting.
$ws.Range("H2").Value = ""
$ws.Range("H2").ClearFormats()

$ws.Range("I2").Value = "δοκιμασία"
$ws.Range("I2").ClearFormats()

# Rows 3-4: de_DE (H3/H4) previously empty, now filled in with the fr text
# that was reused for this column.
$ws.Range("H3").Value = "test 1 fr"
$ws.Range("H4").Value = "test 1 fr"
